$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the values between row 2 and row 3 for columns D, M, N, O, P, Q, S, T

# Column D (Fecha)
$ws.Range("D2").Value = 44855
$ws.Range("D3").Value = 44874

# Column M (Volumen)
$ws.Range("M2").Value = 25
$ws.Range("M3").Value = 67

# Column N (Precio minimo)
$ws.Range("N2").Value = 15000
$ws.Range("N3").Value = 16000

# Column O (Precio maximo)
$ws.Range("O2").Value = 15000
$ws.Range("O3").Value = 16000

# Column P (Precio promedio ponderado)
$ws.Range("P2").Value = 15000
$ws.Range("P3").Value = 16000

# Column Q (Unidad de comercializacion)
$ws.Range("Q2").Value = "`$/bandeja 5 kilos"
$ws.Range("Q3").Value = "`$/bandeja 10 kilos"

# Column S (Precio $/Kg)
$ws.Range("S2").Value = 3000
$ws.Range("S3").Value = 1600

# Column T (Kg / unidad)
$ws.Range("T2").Value = 5
$ws.Range("T3").Value = 10
